$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = -21.42300000000001
    32  = -21.19089999999999
    36  = -19.989
    38  = -19.85659999999998
    46  = -21.89370000000001
    54  = -22.1818
    55  = -22.27830000000001
    67  = -21.40049999999998
    69  = -21.53139999999998
    72  = -21.67859999999999
    91  = -20.62769999999999
    99  = -21.8373
    104 = -21.17279999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}

$wb.Save()
